$d = $word.ActiveDocument

$replacements = @(
    @("187×8=", "184×3="),
    @("919×9=", "312×7="),
    @("241×2=", "596×7="),
    @("806×4=", "447×5="),
    @("319×7=", "794×6="),
    @("530×2=", "593×5="),
    @("962×9=", "316×4="),
    @("882×9=", "583×7="),
    @("296×6=", "850×2="),
    @("724×9=", "881×2="),
    @("741×2=", "280×6="),
    @("141×6=", "525×6="),
    @("436×4=", "661×6="),
    @("602×3=", "683×6="),
    @("360×9=", "851×9="),
    @("830×7=", "224×6="),
    @("805×4=", "347×5="),
    @("705×7=", "369×2="),
    @("563×4=", "726×3="),
    @("181×5=", "869×5="),
    @("173×2=", "107×4="),
    @("280×7=", "577×6="),
    @("581×6=", "978×5="),
    @("456×2=", "385×2="),
    @("709×6=", "699×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
